$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-derived-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Move the "ele-1/ext-1" constraint off of the top-level Extension row (row 2)
# and attach it to the Extension.extension row (row 4) instead.
# (Set via a leading apostrophe + format-only paste from a neighboring blank
# cell so the result is a genuine empty text value - not a fully blank cell -
# while keeping the original cell style/number format untouched.)
$elements.Range("AI2").Value = "'"
$elements.Range("AH2").Copy()
$elements.Range("AI2").PasteSpecial(-4122)
$elements.Range("AI4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# The canonical/base URL also appears as the Fixed Value of Extension.url (row 5); keep it in sync.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-derived-indicator"
